$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A216").Value = 44776
$ws.Range("B216").Value = 'KA 01 MG 3419'
$ws.Range("C216").Value = 'I20'
$ws.Range("D216").Value = 'TIE MEMBER                INSU'
$ws.Range("E216").Value = 'WORK IN PROGRESS'

$ws.Range("A217").Value = 44777
$ws.Range("B217").Value = 'KA 53 MA 4313'
$ws.Range("C217").Value = 'RITZ'
$ws.Range("D217").Value = 'GENERAL CHECKUP         WW'
$ws.Range("E217").Value = 'WORK IN PROGRESS'

$ws.Range("A218").Value = 44777
$ws.Range("B218").Value = 'MH 14 CK 5854'
$ws.Range("C218").Value = 'I10'
$ws.Range("D218").Value = 'GENERAL CHECKUP'
$ws.Range("E218").Value = 'WORK DONE DELIVERED'
$ws.Range("F218").Value = 2589
$ws.Range("G218").Value = 'G PAY'

$ws.Range("A219").Value = 44777
$ws.Range("B219").Value = 'KA03NE7311'
$ws.Range("C219").Value = 'CRYSTA'
$ws.Range("D219").Value = 'BODY SHOP'
$ws.Range("E219").Value = 'WORK DONE'
$ws.Range("F219").Value = 18880
$ws.Range("G219").Value = '  INSURANCE'

$ws.Range("A220").Value = 44777
$ws.Range("B220").Value = 'KA03MU7732'
$ws.Range("C220").Value = 'NISSON TERRANO'
$ws.Range("D220").Value = 'PMS'
$ws.Range("E220").Value = 'WORK DONE DELIVERED'
$ws.Range("F220").Value = 24819

$ws.Range("A221").Value = 44777
$ws.Range("B221").Value = 'KA03MT2522'
$ws.Range("C221").Value = 'ETIOS LIVA'
$ws.Range("D221").Value = 'BODY SHOP'
$ws.Range("E221").Value = 'WORK DONE'
$ws.Range("F221").Value = 14641
$ws.Range("G221").Value = '  INSURANCE'

$ws.Range("A222").Value = 44777
$ws.Range("B222").Value = 'KA01MN9272'
$ws.Range("C222").Value = 'SWIFT DZIRE'
$ws.Range("D222").Value = 'BODY SHOP'
$ws.Range("E222").Value = 'WORK DONE '
$ws.Range("F222").Value = 7375
$ws.Range("G222").Value = '  INSURANCE'

$ws.Range("A223").Value = 44777
$ws.Range("B223").Value = 'LK648697'
$ws.Range("C223").Value = 'BEAT'
$ws.Range("D223").Value = 'PMS'
$ws.Range("E223").Value = 'WORK DONE DELIVERED'
$ws.Range("F223").Value = 7987
$ws.Range("G223").Value = 'G PAY'

$ws.Range("A224").Value = 44781
$ws.Range("B224").Value = 'KA53MC0494'
$ws.Range("C224").Value = 'LINEA'
$ws.Range("D224").Value = 'RR'
$ws.Range("E224").Value = 'WORK DONE DELIVERED'
$ws.Range("F224").Value = 6250
$ws.Range("G224").Value = 'G PAY'

$ws.Range("A225").Value = 44781
$ws.Range("B225").Value = 'KA03MR907'
$ws.Range("C225").Value = 'POLO'
$ws.Range("D225").Value = 'GENERAL CHECKUP'
$ws.Range("E225").Value = 'WORK DONE DELIVERED'
$ws.Range("F225").Value = 3016

$ws.Range("A226").Value = 44781
$ws.Range("B226").Value = 'KA04MS9197'
$ws.Range("C226").Value = 'RAPID'
$ws.Range("D226").Value = 'ABS SENSOR'
$ws.Range("E226").Value = 'WORK DONE DELIVERED'
$ws.Range("F226").Value = 5562

$ws.Range("A227").Value = 44781
$ws.Range("B227").Value = 'KA53MA2760'
$ws.Range("C227").Value = 'BEAT'
$ws.Range("D227").Value = 'RR'
$ws.Range("E227").Value = 'WORK DONE DELIVERED'
$ws.Range("F227").Value = 8584
$ws.Range("G227").Value = 'P PAY'

$ws.Range("A228").Value = 44781
$ws.Range("B228").Value = 'KA53MA4313'
$ws.Range("C228").Value = 'RITZ'
$ws.Range("D228").Value = 'RR'
$ws.Range("E228").Value = 'WORK DONE DELIVERED'
$ws.Range("F228").Value = 5637

$ws.Range("A229").Value = 44781
$ws.Range("B229").Value = 'KA04MR6014'
$ws.Range("C229").Value = 'VENTO'
$ws.Range("D229").Value = 'WIPER BLADE CHANGE'
$ws.Range("E229").Value = 'WORK DONE DELIVERED'
$ws.Range("F229").Value = 1040

$ws.Range("A230").Value = 44781
$ws.Range("B230").Value = 'KA03MZ9550'
$ws.Range("C230").Value = 'ECOSPORT'
$ws.Range("D230").Value = 'PMS                                      WW'
$ws.Range("E230").Value = 'WORK DONE DELIVERED'
$ws.Range("F230").Value = 25265
$ws.Range("G230").Value = '      CREDIT'

$ws.Range("A231").Value = 44781
$ws.Range("B231").Value = 'KA51MB4552'
$ws.Range("C231").Value = 'POLO'
$ws.Range("D231").Value = 'RR'
$ws.Range("E231").Value = 'WORK DONE DELIVERED'
$ws.Range("F231").Value = 23380

$ws.Range("A232").Value = 44781
$ws.Range("B232").Value = 'PY01CC6847'
$ws.Range("C232").Value = 'SWIFT'
$ws.Range("D232").Value = 'LED LIGHT CHANGE'
$ws.Range("E232").Value = 'WORK DONE DELIVERED'
$ws.Range("F232").Value = 5500
$ws.Range("G232").Value = 'G PAY'

$ws.Range("A233").Value = 44781
$ws.Range("B233").Value = 'DL5CN9218'
$ws.Range("C233").Value = 'BALENO'
$ws.Range("D233").Value = 'PMS                                      WW'
$ws.Range("E233").Value = 'WORK DONE DELIVERED'
$ws.Range("F233").Value = 4243
$ws.Range("G233").Value = 'CREDIT'

$ws.Range("A234").Value = 44781
$ws.Range("B234").Value = 'KA53MQ8210'
$ws.Range("C234").Value = 'I20'
$ws.Range("D234").Value = 'BODY SHOP'
$ws.Range("E234").Value = 'WORK DONE DELIVERED'
$ws.Range("F234").Value = 162681

$ws.Range("A235").Value = 44781
$ws.Range("B235").Value = 'KA53MH7954'
$ws.Range("C235").Value = 'NEXON'
$ws.Range("D235").Value = 'PMS'
$ws.Range("E235").Value = 'WORK DONE DELIVERED'
$ws.Range("F235").Value = 3985

$ws.Range("A236").Value = 44781
$ws.Range("B236").Value = 'KA53Z4971'
$ws.Range("C236").Value = 'INDICA VISTA'
$ws.Range("D236").Value = 'RR'
$ws.Range("E236").Value = 'WORK DONE DELIVERED'
$ws.Range("F236").Value = 1636
$ws.Range("B236").NumberFormat = $ws.Range("A236").NumberFormat

$ws.Range("A237").Value = 44781
$ws.Range("B237").Value = 'KA01MR8095'
$ws.Range("C237").Value = 'STROME'
$ws.Range("D237").Value = 'RR'
$ws.Range("E237").Value = 'WORK DONE DELIVERED'
$ws.Range("F237").Value = 7316

$ws.Range("A238").Value = 44782
$ws.Range("B238").Value = 'MH10AN9645'
$ws.Range("C238").Value = 'FABIA'
$ws.Range("D238").Value = 'BODY SHOP'
$ws.Range("E238").Value = 'WORK DONE'
$ws.Range("F238").Value = 14333
$ws.Range("G238").Value = '  INSURANCE'

$ws.Range("A239").Value = 44782
$ws.Range("B239").Value = 'KA51MD0287'
$ws.Range("C239").Value = 'INNOVA'
$ws.Range("D239").Value = 'RR'
$ws.Range("E239").Value = 'WORK DONE DELIVERED'
$ws.Range("F239").Value = 49087
$ws.Range("G239").Value = 'CARD'

$ws.Range("A240").Value = 44782
$ws.Range("B240").Value = 'KA03MM7229'
$ws.Range("C240").Value = 'I10'
$ws.Range("D240").Value = 'PMS'
$ws.Range("E240").Value = 'WORK DONE DELIVERED'
$ws.Range("F240").Value = 2733
$ws.Range("G240").Value = 'CREDIT'

$ws.Range("A241").Value = 44783
$ws.Range("B241").Value = 'KA04MK7854'
$ws.Range("C241").Value = 'VENTO'
$ws.Range("D241").Value = 'PMS'
$ws.Range("E241").Value = 'WORK DONE DELIVERED'
$ws.Range("F241").Value = 22576
$ws.Range("G241").Value = 'CARD'

$ws.Range("A242").Value = 44783
$ws.Range("B242").Value = 'KA03MJ4271'
$ws.Range("C242").Value = 'FIESTA'
$ws.Range("D242").Value = 'GENERAL CHECKUP'
$ws.Range("E242").Value = 'WORK DONE DELIVERED'
$ws.Range("F242").Value = 944
$ws.Range("G242").Value = 'P PAY'

$ws.Range("A243").Value = 44783
$ws.Range("B243").Value = 'PY 01 BL 1480 '
$ws.Range("C243").Value = 'XYLO'
$ws.Range("D243").Value = 'RR'
$ws.Range("E243").Value = 'WORK DONE DELIVERED'
$ws.Range("F243").Value = 9180
$ws.Range("G243").Value = 'G PAY'

$ws.Range("A244").Value = 44783
$ws.Range("B244").Value = 'KA03ND9387'
$ws.Range("C244").Value = 'DUSTER'
$ws.Range("D244").Value = 'PMS                                      WW'
$ws.Range("E244").Value = 'WORK DONE DELIVERED'
$ws.Range("F244").Value = 5010
$ws.Range("G244").Value = 'CREDIT'

$ws.Range("A245").Value = 44784
$ws.Range("B245").Value = 'KA03MQ0016'
$ws.Range("C245").Value = 'POLO'
$ws.Range("D245").Value = 'PMS                                      WW'
$ws.Range("E245").Value = 'WORK DONE DELIVERED'
$ws.Range("F245").Value = 11026
$ws.Range("G245").Value = 'CREDIT'

$ws.Range("A246").Value = 44784
$ws.Range("B246").Value = 'KA51MB2926'
$ws.Range("C246").Value = 'BEAT'
$ws.Range("D246").Value = 'PMS'
$ws.Range("E246").Value = 'WORK DONE DELIVERED'
$ws.Range("F246").Value = 18659
$ws.Range("G246").Value = 'CARD'

$ws.Range("A247").Value = 44784
$ws.Range("B247").Value = 'KA01MG7422'
$ws.Range("C247").Value = 'SUPERB'
$ws.Range("D247").Value = 'SCANNING'
$ws.Range("E247").Value = 'WORK DONE DELIVERED'
$ws.Range("F247").Value = 3875
$ws.Range("G247").Value = 'CARD'

$ws.Range("A248").Value = 44784
$ws.Range("B248").Value = 'KA05MS4170'
$ws.Range("C248").Value = 'CELERIO'
$ws.Range("D248").Value = 'HORN CHANGE'
$ws.Range("E248").Value = 'WORK DONE DELIVERED'
$ws.Range("F248").Value = 1340
$ws.Range("G248").Value = 'G PAY'

# Column A width adjustment (best achievable precision in this engine)
$ws.Columns.Item(1).ColumnWidth = 11.86

# Update selection to final cell, matching the author's last edit position
$ws.Range("H248").Select()